$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '28.243.02'
$ws.Range('E2').Value = '  -2.61%  '
$ws.Range('D3').Value = '1.868.02'
$ws.Range('D4').Value = '1.004'
$ws.Range('E4').Value = '  +0.12%  '
$ws.Range('D5').Value = '318.59'
$ws.Range('E5').Value = '  -1.79%  '
$ws.Range('E6').Value = '  +0.21%  '
$ws.Range('D7').Value = '0.4395'
$ws.Range('E7').Value = '  -4.24%  '
$ws.Range('D8').Value = '0.3691'
$ws.Range('E8').Value = '  -3.48%  '
$ws.Range('D9').Value = '0.07498'
$ws.Range('E9').Value = '  -2.76%  '
$ws.Range('D10').Value = '0.9352'
$ws.Range('E10').Value = '  -4.55%  '
$ws.Range('D11').Value = '21.33'
$ws.Range('D12').Value = '1.902.65'
$ws.Range('E12').Value = '  +0.73%  '
$ws.Range('B13').Value = 'Chainlink'
$ws.Range('C13').Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$ws.Range('D13').Value = '6.711'
$ws.Range('E13').Value = '  -3.21%  '
$ws.Range('B14').Value = 'Polkadot'
$ws.Range('C14').Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range('D14').Value = '5.471'
$ws.Range('E14').Value = '  -3.58%  '
$ws.Range('D15').Value = '0.06899'
$ws.Range('E15').Value = '  -1.90%  '
$ws.Range('E16').Value = '  +0.20%  '
$ws.Range('D17').Value = '82.00'
$ws.Range('E17').Value = '  -2.21%  '
$ws.Range('D18').Value = '0.000009036'
$ws.Range('E18').Value = '  -4.44%  '
$ws.Range('E19').Value = '  +0.24%  '
$ws.Range('D20').Value = '15.92'
$ws.Range('E20').Value = '  -4.75%  '
$ws.Range('D21').Value = '28.221.51'
$ws.Range('D22').Value = '5.121'
$ws.Range('E22').Value = '  -3.66%  '
$ws.Range('D23').Value = '10.80'
$ws.Range('E23').Value = '  -0.83%  '
$ws.Range('D24').Value = '2.139.01'
$ws.Range('E24').Value = '  +1.03%  '
$ws.Range('E25').Value = '  -3.30%  '
$ws.Range('D26').Value = '155.30'
$ws.Range('E26').Value = '  -1.93%  '
$ws.Range('D27').Value = '18.41'
$ws.Range('E27').Value = '  -3.53%  '
$ws.Range('D28').Value = '5.320'
$ws.Range('E28').Value = '  -6.04%  '
$ws.Range('D29').Value = '113.41'
$ws.Range('E29').Value = '  -3.59%  '
$ws.Range('D30').Value = '1.727'
$ws.Range('E30').Value = '  -6.66%  '
$ws.Range('D31').Value = '0.09026'
$ws.Range('E31').Value = '  -2.78%  '
$ws.Range('D32').Value = '0.7958'
$ws.Range('E32').Value = '  -7.89%  '
$ws.Range('D33').Value = '4.847'
$ws.Range('E33').Value = '  -4.36%  '
$ws.Range('D34').Value = '1.174'
$ws.Range('E34').Value = '  -5.98%  '
$ws.Range('D35').Value = '2.937'
$ws.Range('E35').Value = '  -2.85%  '
$ws.Range('E36').Value = '  +0.22%  '
$ws.Range('D37').Value = '1.128'
$ws.Range('E37').Value = '  -2.35%  '
$ws.Range('D38').Value = '0.05447'
$ws.Range('E38').Value = '  -5.17%  '
$ws.Range('D39').Value = '0.01968'
$ws.Range('E39').Value = '  -3.69%  '
$ws.Range('E40').Value = '  +3.32%  '
$ws.Range('D41').Value = '0.5265'
$ws.Range('E41').Value = '  -4.43%  '
$ws.Range('D42').Value = '7.058'
$ws.Range('E42').Value = '  -4.85%  '
$ws.Range('E43').Value = '  -4.32%  '
$ws.Range('D44').Value = '8.720'
$ws.Range('E44').Value = '  -6.52%  '
$ws.Range('D45').Value = '0.06745'
$ws.Range('E45').Value = '  -1.32%  '
$ws.Range('D46').Value = '0.4874'
$ws.Range('E46').Value = '  -5.91%  '
$ws.Range('D47').Value = '10.53'
$ws.Range('E47').Value = '  -6.42%  '
$ws.Range('D48').Value = '107.01'
$ws.Range('E48').Value = '  -3.55%  '
$ws.Range('D49').Value = '1.918'
$ws.Range('E49').Value = '  -6.59%  '
$ws.Range('D50').Value = '1.003'
$ws.Range('E50').Value = '  +0.17%  '
$ws.Range('D51').Value = '1.674'
$ws.Range('E51').Value = '  -6.08%  '
